$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "(""Black Sun's Zenith"", ['{X}{B}{B}', 'Sorcery', 'Put X -1/-1 counters on each creature. Shuffle Black Sun" + [char]8217 + "s Zenith into its owner" + [char]8217 + "s library.'])"
$ws.Range("A3").Value = "('Glissa, the Traitor', ['{B}{G}{G}', 'Legendary Creature " + [char]8212 + " Zombie Elf', 'First strike, deathtouch', 'Whenever a creature an opponent controls dies, you may return target artifact card from your graveyard to your hand.', '3/3'])"
$ws.Range("A4").Value = "('Hero of Bladehold', ['{2}{W}{W}', 'Creature " + [char]8212 + " Human Knight', 'Battle cry (Whenever this creature attacks, each other attacking creature gets +1/+0 until end of turn.)', 'Whenever Hero of Bladehold attacks, create two 1/1 white Soldier creature tokens that are tapped and attacking.', '3/4'])"
$ws.Range("A5").Value = "('Mirran Crusader', ['{1}{W}{W}', 'Creature " + [char]8212 + " Human Knight', 'Double strike, protection from black and from green', '2/2'])"
$ws.Range("A6").Value = "('Thopter Assembly', ['{6}', 'Artifact Creature " + [char]8212 + " Thopter', 'Flying', 'At the beginning of your upkeep, if you control no Thopters other than Thopter Assembly, return Thopter Assembly to its owner" + [char]8217 + "s hand and create five 1/1 colorless Thopter artifact creature tokens with flying.', '5/5'])"
$ws.Range("A7").Value = "('Treasure Mage', ['{2}{U}', 'Creature " + [char]8212 + " Human Wizard', 'When Treasure Mage enters the battlefield, you may search your library for an artifact card with converted mana cost 6 or greater, reveal it, put it into your hand, then shuffle your library.', '2/2'])"

# Delete rows 8 through 33 (old leftover data)
$ws.Range("A8:A33").EntireRow.Delete()
